$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Lif"
$ws.Range("C2").Value = "Il6st"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.904627
$ws.Range("H2").Value = 2.713881
$ws.Range("I2").Value = 0.2670571014571191
$ws.Range("J2").Value = 0.2670571014571191
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 37.42645866666667
$ws.Range("N2").Value = 112.279376
$ws.Range("O2").Value = 0.2415534622699011
$ws.Range("P2").Value = 0.2415534622699011
$ws.Range("Q2").Value = 33.85698502425067
$ws.Range("R2").Value = 304.712865218256
$ws.Range("S2").Value = 0.06450856748073137
$ws.Range("T2").Value = 0.06450856748073137

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Lif"
$ws.Range("C3").Value = "Il6st"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.904627
$ws.Range("H3").Value = 2.713881
$ws.Range("I3").Value = 0.2670571014571191
$ws.Range("J3").Value = 0.2670571014571191
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 94.96115633333334
$ws.Range("N3").Value = 284.883469
$ws.Range("O3").Value = 0.6128871635375853
$ws.Range("P3").Value = 0.6128871635375853
$ws.Range("Q3").Value = 85.90442597035432
$ws.Range("R3").Value = 773.1398337331889
$ws.Range("S3").Value = 0.1636758694146229
$ws.Range("T3").Value = 0.1636758694146229

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Lif"
$ws.Range("C4").Value = "Il6st"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.904627
$ws.Range("H4").Value = 2.713881
$ws.Range("I4").Value = 0.2670571014571191
$ws.Range("J4").Value = 0.2670571014571191
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 22.553069
$ws.Range("N4").Value = 67.659207
$ws.Range("O4").Value = 0.1455593741925136
$ws.Range("P4").Value = 0.1455593741925136
$ws.Range("Q4").Value = 20.402115150263
$ws.Range("R4").Value = 183.619036352367
$ws.Range("S4").Value = 0.03887266456176486
$ws.Range("T4").Value = 0.03887266456176486

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Lif"
$ws.Range("C5").Value = "Il6st"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.482764666666667
$ws.Range("H5").Value = 7.448294000000001
$ws.Range("I5").Value = 0.7329428985428809
$ws.Range("J5").Value = 0.7329428985428807
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 37.42645866666667
$ws.Range("N5").Value = 112.279376
$ws.Range("O5").Value = 0.2415534622699011
$ws.Range("P5").Value = 0.2415534622699011
$ws.Range("Q5").Value = 92.92108917606045
$ws.Range("R5").Value = 836.2898025845441
$ws.Range("S5").Value = 0.1770448947891697
$ws.Range("T5").Value = 0.1770448947891697

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Lif"
$ws.Range("C6").Value = "Il6st"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.482764666666667
$ws.Range("H6").Value = 7.448294000000001
$ws.Range("I6").Value = 0.7329428985428809
$ws.Range("J6").Value = 0.7329428985428807
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 94.96115633333334
$ws.Range("N6").Value = 284.883469
$ws.Range("O6").Value = 0.6128871635375853
$ws.Range("P6").Value = 0.6128871635375853
$ws.Range("Q6").Value = 235.7662036502096
$ws.Range("R6").Value = 2121.895832851886
$ws.Range("S6").Value = 0.4492112941229624
$ws.Range("T6").Value = 0.4492112941229623

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Lif"
$ws.Range("C7").Value = "Il6st"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.482764666666667
$ws.Range("H7").Value = 7.448294000000001
$ws.Range("I7").Value = 0.7329428985428809
$ws.Range("J7").Value = 0.7329428985428807
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 22.553069
$ws.Range("N7").Value = 67.659207
$ws.Range("O7").Value = 0.1455593741925136
$ws.Range("P7").Value = 0.1455593741925136
$ws.Range("Q7").Value = 55.99396283809533
$ws.Range("R7").Value = 503.945665542858
$ws.Range("S7").Value = 0.1066867096307487
$ws.Range("T7").Value = 0.1066867096307487
